$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.368.05"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.837.84"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.58"
$ws.Range("E5").Value = "  +5.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.24"
$ws.Range("E6").Value = "  -5.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.603"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.704"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -6.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("E11").Value = "  -7.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.22"
$ws.Range("E12").Value = "  -4.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.22"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "4.441.52"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.60"
$ws.Range("E15").Value = "  +7.98%  "
$ws.Range("D16").Value = "3.834.50"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.11"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.20"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "68.385.88"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "415.22"
$ws.Range("E21").Value = "  -4.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.43"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.90"
$ws.Range("E23").Value = "  -5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.13"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +5.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.48"
$ws.Range("E26").Value = "  -7.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.19"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.05"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "669.06"
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.124"
$ws.Range("E31").Value = "  -6.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.09"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.441"
$ws.Range("E35").Value = "  -10.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.51"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").Value = "0.0₃0833"
$ws.Range("E37").Value = "  -6.17%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.147"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +9.63%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0470"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("E44").Value = "  -6.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.138"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.71"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000269"
$ws.Range("E49").Value = "  +11.70%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.662.85"
$ws.Range("E51").Value = "  +9.45%  "
